# Update quizvragen via Admin
# DC sheet: insert a new "q1" (Ohmse wet) question as row 2, pushing the
# existing "q3" (Vermogen) question down to row 3, and append a brand new
# draft question ("Is dit een goede nieuwe vraag??") as row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure we're editing the "DC" sheet (it's already the active sheet in
# this workbook, but be explicit/defensive).
if ($ws.Name -ne "DC") {
    $ws = $wb.Worksheets.Item("DC")
}

# --- Insert a fresh row above the current row 2 ("q3") -------------------
$ws.Rows.Item(2).Insert()
# Inserting copies formatting down from the row above (the bold header);
# the source data has no explicit styling on data rows, so strip it back
# off to keep the new row plain, matching the rest of the sheet.
$ws.Rows.Item(2).ClearFormats()

# --- Row 2: new question "q1" - Ohmse wet --------------------------------
$ws.Cells.Item(2, 1).Value  = "q1"
$ws.Cells.Item(2, 2).Value  = "mc"
$ws.Cells.Item(2, 3).Value  = "Ohmse wet"
$ws.Cells.Item(2, 4).Value  = "Wat is de juiste formule voor de stroom I?"
$ws.Cells.Item(2, 5).Value  = "['I = U/R', 'U = I*R', 'R = U/I']"
$ws.Cells.Item(2, 6).Value  = 0
$ws.Cells.Item(2, 7).Value  = "Volgens de wet van Ohm geldt: I = U / R."
$ws.Cells.Item(2, 8).Value  = "assets/ohm_schema.png"
$ws.Cells.Item(2, 9).Value  = "I = \frac{U}{R}"
$ws.Cells.Item(2, 10).Value = "['DC','basis']"
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/DC_q0_1763135501.png"

# Row 3 ("q3" - Vermogen) already shifted down automatically by the Insert
# above; its values (and the existing empty image_url cell) are untouched.

# --- Row 4: brand-new draft question appended at the end -----------------
$ws.Cells.Item(4, 2).Value  = "mc"
$ws.Cells.Item(4, 4).Value  = "Is dit een goede nieuwe vraag??"
$ws.Cells.Item(4, 5).Value  = "['A. test 1', 'B. Test 2', 'C. Test 3']"
$ws.Cells.Item(4, 6).Value  = 0
$ws.Cells.Item(4, 12).Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/DC_new_1763469796.jpg"
